# "Avancement sur le jeu"
# Fill in a new "Début" (start) entry for row 27 of the "Activités" sheet
# (Date + start time), matching the ongoing work started 2021-03-19 10:37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activités")

# Row 27: Date (A27) and Début / start time (B27)
# (19 mars 2021 / 10:37 as Excel date-time serials, matching the existing
# column formatting already applied to these cells)
$ws.Cells.Item(27, 1).Value = 44274
$ws.Cells.Item(27, 2).Value = 0.44236111111111115

# Recalculate the workbook so the volatile "Temps" formulas (NOW()-based,
# rows 26 & 27) pick up the current timestamp.
$excel.CalculateFullRebuild()

# Update the active selection to reflect where the user left off editing.
$ws.Activate()
$ws.Range("E27").Select()
